# feat: new column (%) in report
#
# 1. qc sheet: report_date (column F) values bumped to the new run timestamp.
# 2. emu_long sheet: new column F "% total" = "abundance total" (col E) * 100,
#    expressed as a percentage value.

$wb = $excel.ActiveWorkbook

# ---- qc sheet: refresh the report_date timestamp (column F, rows 2-6) ----
$qc = $wb.Worksheets.Item("qc")
$qc.Range("F2:F6").Value = 45642.47677915011

# ---- emu_long sheet: add "% total" column ----
$long = $wb.Worksheets.Item("emu_long")

$long.Range("F1").Value = "% total"

$pctTotals = @{
    2  = 100
    3  = 0
    4  = 100
    5  = 100
    6  = 0
    7  = 100
    8  = 100
    9  = 0
    10 = 100
    11 = 88.48560700876095
    12 = 6.774023274569065
    13 = 4.740369716669985
    14 = 0
    15 = 100
    16 = 44.44444444444444
    17 = 30.81387123397617
    18 = 24.74168432157938
    19 = 0
    20 = 100
}

foreach ($row in $pctTotals.Keys) {
    $long.Cells.Item($row, 6).Value = $pctTotals[$row]
}
